$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = "ba"
$ws.Range("J19").Value = "Appreciation"
$ws.Range("I22").Value = "ba"
$ws.Range("J22").Value = "Appreciation"
$ws.Range("I23").Value = "sv"
$ws.Range("J23").Value = "Statement-opinion"
$ws.Range("I26").Value = "ba"
$ws.Range("J26").Value = "Appreciation"
$ws.Range("I28").Value = "b"
$ws.Range("J28").Value = "Acknowledge (Backchannel)"
$ws.Range("I32").Value = "sd"
$ws.Range("J32").Value = "Statement-non-opinion"
$ws.Range("I40").Value = "aa"
$ws.Range("J40").Value = "Agree/Accept"
$ws.Range("I45").Value = "aa"
$ws.Range("J45").Value = "Agree/Accept"
$ws.Range("I46").Value = "aa"
$ws.Range("J46").Value = "Agree/Accept"
$ws.Range("I47").Value = "b"
$ws.Range("J47").Value = "Acknowledge (Backchannel)"
$ws.Range("I50").Value = "aa"
$ws.Range("J50").Value = "Agree/Accept"
$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I60").Value = "sd"
$ws.Range("J60").Value = "Statement-non-opinion"
$ws.Range("I62").Value = "sv"
$ws.Range("J62").Value = "Statement-opinion"
$ws.Range("I64").Value = "aa"
$ws.Range("J64").Value = "Agree/Accept"
$ws.Range("I65").Value = "sd"
$ws.Range("J65").Value = "Statement-non-opinion"
$ws.Range("I67").Value = "b"
$ws.Range("J67").Value = "Acknowledge (Backchannel)"
$ws.Range("I71").Value = "sv"
$ws.Range("J71").Value = "Statement-opinion"
$ws.Range("I74").Value = "sd"
$ws.Range("J74").Value = "Statement-non-opinion"
$ws.Range("I75").Value = "sv"
$ws.Range("J75").Value = "Statement-opinion"
$ws.Range("I84").Value = "sv"
$ws.Range("J84").Value = "Statement-opinion"
$ws.Range("I85").Value = "sv"
$ws.Range("J85").Value = "Statement-opinion"
$ws.Range("I91").Value = "b"
$ws.Range("J91").Value = "Acknowledge (Backchannel)"
$ws.Range("I94").Value = "ba"
$ws.Range("J94").Value = "Appreciation"
$ws.Range("I100").Value = "sd"
$ws.Range("J100").Value = "Statement-non-opinion"
$ws.Range("I101").Value = "sd"
$ws.Range("J101").Value = "Statement-non-opinion"
$ws.Range("I107").Value = "sd"
$ws.Range("J107").Value = "Statement-non-opinion"
$ws.Range("I114").Value = "sd"
$ws.Range("J114").Value = "Statement-non-opinion"
$ws.Range("I129").Value = "sd"
$ws.Range("J129").Value = "Statement-non-opinion"
$ws.Range("I165").Value = "ba"
$ws.Range("J165").Value = "Appreciation"
$ws.Range("I171").Value = "ba"
$ws.Range("J171").Value = "Appreciation"
$ws.Range("I188").Value = "%"
$ws.Range("J188").Value = "Uninterpretable"
$ws.Range("I196").Value = "sd"
$ws.Range("J196").Value = "Statement-non-opinion"
$ws.Range("I201").Value = "sd"
$ws.Range("J201").Value = "Statement-non-opinion"
$ws.Range("I203").Value = "sd"
$ws.Range("J203").Value = "Statement-non-opinion"
$ws.Range("I220").Value = "sd"
$ws.Range("J220").Value = "Statement-non-opinion"
$ws.Range("I236").Value = "aa"
$ws.Range("J236").Value = "Agree/Accept"
$ws.Range("I241").Value = "ba"
$ws.Range("J241").Value = "Appreciation"
$ws.Range("I243").Value = "ba"
$ws.Range("J243").Value = "Appreciation"
$ws.Range("I244").Value = "sv"
$ws.Range("J244").Value = "Statement-opinion"
$ws.Range("I247").Value = "ba"
$ws.Range("J247").Value = "Appreciation"
$ws.Range("I252").Value = "sd"
$ws.Range("J252").Value = "Statement-non-opinion"
$ws.Range("I257").Value = "sv"
$ws.Range("J257").Value = "Statement-opinion"
$ws.Range("I260").Value = "sd"
$ws.Range("J260").Value = "Statement-non-opinion"
$ws.Range("I276").Value = "ba"
$ws.Range("J276").Value = "Appreciation"
$ws.Range("I279").Value = "sd"
$ws.Range("J279").Value = "Statement-non-opinion"
$ws.Range("I282").Value = "sd"
$ws.Range("J282").Value = "Statement-non-opinion"
$ws.Range("I283").Value = "ba"
$ws.Range("J283").Value = "Appreciation"
$ws.Range("I285").Value = "b"
$ws.Range("J285").Value = "Acknowledge (Backchannel)"
$ws.Range("I287").Value = "b"
$ws.Range("J287").Value = "Acknowledge (Backchannel)"
$ws.Range("I290").Value = "b"
$ws.Range("J290").Value = "Acknowledge (Backchannel)"
$ws.Range("I294").Value = "sd"
$ws.Range("J294").Value = "Statement-non-opinion"
$ws.Range("I299").Value = "sv"
$ws.Range("J299").Value = "Statement-opinion"
$ws.Range("I301").Value = "b"
$ws.Range("J301").Value = "Acknowledge (Backchannel)"
$ws.Range("I304").Value = "ba"
$ws.Range("J304").Value = "Appreciation"
$ws.Range("I307").Value = "sd"
$ws.Range("J307").Value = "Statement-non-opinion"
$ws.Range("I308").Value = "sd"
$ws.Range("J308").Value = "Statement-non-opinion"
$ws.Range("I310").Value = "sd"
$ws.Range("J310").Value = "Statement-non-opinion"
$ws.Range("I315").Value = "sd"
$ws.Range("J315").Value = "Statement-non-opinion"
$ws.Range("I317").Value = "aa"
$ws.Range("J317").Value = "Agree/Accept"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I323").Value = "aa"
$ws.Range("J323").Value = "Agree/Accept"
$ws.Range("I332").Value = "ba"
$ws.Range("J332").Value = "Appreciation"
$ws.Range("I334").Value = "ba"
$ws.Range("J334").Value = "Appreciation"
$ws.Range("I340").Value = "ba"
$ws.Range("J340").Value = "Appreciation"
$ws.Range("I358").Value = "aa"
$ws.Range("J358").Value = "Agree/Accept"
$ws.Range("I359").Value = "%"
$ws.Range("J359").Value = "Uninterpretable"
$ws.Range("I368").Value = "ba"
$ws.Range("J368").Value = "Appreciation"
$ws.Range("I369").Value = "b"
$ws.Range("J369").Value = "Acknowledge (Backchannel)"
$ws.Range("I371").Value = "sv"
$ws.Range("J371").Value = "Statement-opinion"
$ws.Range("I372").Value = "sv"
$ws.Range("J372").Value = "Statement-opinion"
$ws.Range("I374").Value = "sv"
$ws.Range("J374").Value = "Statement-opinion"
$ws.Range("I375").Value = "sv"
$ws.Range("J375").Value = "Statement-opinion"
$ws.Range("I376").Value = "sv"
$ws.Range("J376").Value = "Statement-opinion"
$ws.Range("I389").Value = "sd"
$ws.Range("J389").Value = "Statement-non-opinion"
$ws.Range("I397").Value = "ba"
$ws.Range("J397").Value = "Appreciation"
$ws.Range("I416").Value = "sd"
$ws.Range("J416").Value = "Statement-non-opinion"
$ws.Range("I425").Value = "aa"
$ws.Range("J425").Value = "Agree/Accept"
$ws.Range("I433").Value = "b"
$ws.Range("J433").Value = "Acknowledge (Backchannel)"
$ws.Range("I451").Value = "sd"
$ws.Range("J451").Value = "Statement-non-opinion"
$ws.Range("I455").Value = "ba"
$ws.Range("J455").Value = "Appreciation"
$ws.Range("I463").Value = "aa"
$ws.Range("J463").Value = "Agree/Accept"
$ws.Range("I465").Value = "sv"
$ws.Range("J465").Value = "Statement-opinion"
$ws.Range("I486").Value = "sv"
$ws.Range("J486").Value = "Statement-opinion"
$ws.Range("I491").Value = "aa"
$ws.Range("J491").Value = "Agree/Accept"
$ws.Range("I494").Value = "ba"
$ws.Range("J494").Value = "Appreciation"
$ws.Range("I496").Value = "aa"
$ws.Range("J496").Value = "Agree/Accept"
$ws.Range("I501").Value = "sv"
$ws.Range("J501").Value = "Statement-opinion"
$ws.Range("I506").Value = "sd"
$ws.Range("J506").Value = "Statement-non-opinion"
$ws.Range("I516").Value = "sd"
$ws.Range("J516").Value = "Statement-non-opinion"
$ws.Range("I521").Value = "ba"
$ws.Range("J521").Value = "Appreciation"
$ws.Range("I522").Value = "b"
$ws.Range("J522").Value = "Acknowledge (Backchannel)"
$ws.Range("I524").Value = "sd"
$ws.Range("J524").Value = "Statement-non-opinion"
$ws.Range("I526").Value = "sv"
$ws.Range("J526").Value = "Statement-opinion"
$ws.Range("I530").Value = "sd"
$ws.Range("J530").Value = "Statement-non-opinion"
$ws.Range("I531").Value = "sd"
$ws.Range("J531").Value = "Statement-non-opinion"
